$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Range("D:E").Insert()

# Copy the number formatting/style from column F (the old column D, now shifted)
# into the two newly inserted columns D:E so the new cells pick up the correct
# date/number styles instead of the generic default column style.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns with the newest two quarters of data
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 1386000
$ws.Cells.Item(8, 5).Value = 1287000
$ws.Cells.Item(9, 4).Value = 844000
$ws.Cells.Item(9, 5).Value = 782000
$ws.Cells.Item(10, 4).Value = 542000
$ws.Cells.Item(10, 5).Value = 505000
$ws.Cells.Item(12, 4).Value = 52000
$ws.Cells.Item(12, 5).Value = 46000
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 3000
$ws.Cells.Item(14, 5).Value = 2000
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(17, 4).Value = 1192000
$ws.Cells.Item(17, 5).Value = 1109000
$ws.Cells.Item(18, 4).Value = 194000
$ws.Cells.Item(18, 5).Value = 178000
$ws.Cells.Item(20, 4).Value = 4000
$ws.Cells.Item(20, 5).Value = 4000
$ws.Cells.Item(21, 4).Value = 264000
$ws.Cells.Item(21, 5).Value = 245000
$ws.Cells.Item(22, 4).Value = 19000
$ws.Cells.Item(22, 5).Value = 21000
$ws.Cells.Item(23, 4).Value = 179000
$ws.Cells.Item(23, 5).Value = 161000
$ws.Cells.Item(24, 4).Value = -53500
$ws.Cells.Item(24, 5).Value = 31000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 232500
$ws.Cells.Item(26, 5).Value = 130000
$ws.Cells.Item(27, 4).Value = 232500
$ws.Cells.Item(27, 5).Value = 130000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = -7500
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = -4000
$ws.Cells.Item(32, 5).Value = -4000
$ws.Cells.Item(33, 4).Value = 225000
$ws.Cells.Item(33, 5).Value = 130000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 225000
$ws.Cells.Item(35, 5).Value = 130000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 296000
$ws.Cells.Item(41, 5).Value = 404000
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 1031000
$ws.Cells.Item(43, 5).Value = 1017000
$ws.Cells.Item(44, 4).Value = 595000
$ws.Cells.Item(44, 5).Value = 622000
$ws.Cells.Item(45, 4).Value = 172000
$ws.Cells.Item(45, 5).Value = 161000
$ws.Cells.Item(46, 4).Value = 2094000
$ws.Cells.Item(46, 5).Value = 2204000
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 656000
$ws.Cells.Item(48, 5).Value = 636000
$ws.Cells.Item(49, 4).Value = 4208000
$ws.Cells.Item(49, 5).Value = 4267000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 264000
$ws.Cells.Item(52, 5).Value = 219000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 7222000
$ws.Cells.Item(54, 5).Value = 7326000
$ws.Cells.Item(57, 4).Value = 586000
$ws.Cells.Item(57, 5).Value = 570000
$ws.Cells.Item(58, 4).Value = 257000
$ws.Cells.Item(58, 5).Value = 342000
$ws.Cells.Item(59, 4).Value = 546000
$ws.Cells.Item(59, 5).Value = 571000
$ws.Cells.Item(60, 4).Value = 1389000
$ws.Cells.Item(60, 5).Value = 1483000
$ws.Cells.Item(61, 4).Value = 2051000
$ws.Cells.Item(61, 5).Value = 2189000
$ws.Cells.Item(62, 4).Value = 1000000
$ws.Cells.Item(62, 5).Value = 1023000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 4454000
$ws.Cells.Item(66, 5).Value = 4709000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 1639000
$ws.Cells.Item(72, 5).Value = 1452000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 2768000
$ws.Cells.Item(76, 5).Value = 2617000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 225000
$ws.Cells.Item(81, 5).Value = 130000
$ws.Cells.Item(83, 4).Value = 66000
$ws.Cells.Item(83, 5).Value = 63000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 198000
$ws.Cells.Item(89, 5).Value = 187000
$ws.Cells.Item(91, 4).Value = -66000
$ws.Cells.Item(91, 5).Value = -60000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -64000
$ws.Cells.Item(94, 5).Value = -42000
$ws.Cells.Item(96, 4).Value = -38000
$ws.Cells.Item(96, 5).Value = -38000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -235000
$ws.Cells.Item(100, 5).Value = -57000
$ws.Cells.Item(101, 4).Value = -7000
$ws.Cells.Item(101, 5).Value = -5000
$ws.Cells.Item(102, 4).Value = -108000
$ws.Cells.Item(102, 5).Value = 83000
